# NYPD 110th Precinct CompStat weekly report refresh:
#   - advance the report's Volume/Number and the covered-week date range
#   - overwrite the weekly/28-day/YTD/2-year crime-statistics grid with the
#     newly collected figures (rows 14-30, 33)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump issue number and the two report dates -------------------
$ws.Range("A8").Value = "Volume 31   Number  38"
$ws.Range("C9").Value = "Report Covering the Week  9/16/2024  Through  9/22/2024"

# --- Phase 1: cells whose fundamental type flips between number and text --
# ("text" cells hold the sheet's "0"/"***.*" placeholders used when a
# count or a percent-change is not meaningful, e.g. division by zero).
# Grab the number format from an unaffected same-row cell that already has
# the target look, then write the real value.

$ws.Range("C18").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("D18").Value = 5

$ws.Range("H18").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = -40

$ws.Range("F22").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C22").Value = 2

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("N22").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("N22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "0"
$ws.Range("D33").Copy()
$ws.Range("C33").PasteSpecial(-4122)

# --- Phase 2: remaining plain numeric updates ------------------------------

$ws.Range("L14").Value = -60
$ws.Range("N14").Value = -91.304347826087

$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -80
$ws.Range("J15").Value = 25
$ws.Range("K15").Value = 12
$ws.Range("M15").Value = 55.555555555555

$ws.Range("C16").Value = 15
$ws.Range("E16").Value = 15.384615384615
$ws.Range("F16").Value = 51
$ws.Range("H16").Value = 10.869565217391
$ws.Range("I16").Value = 385
$ws.Range("J16").Value = 337
$ws.Range("K16").Value = 14.243323442136
$ws.Range("L16").Value = 31.399317406143
$ws.Range("M16").Value = 47.509578544061
$ws.Range("N16").Value = -66.287215411558

$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 18
$ws.Range("E17").Value = -27.777777777777
$ws.Range("F17").Value = 49
$ws.Range("G17").Value = 63
$ws.Range("H17").Value = -22.222222222222
$ws.Range("I17").Value = 605
$ws.Range("J17").Value = 545
$ws.Range("K17").Value = 11.009174311926
$ws.Range("L17").Value = 56.735751295336
$ws.Range("M17").Value = 225.268817204301
$ws.Range("N17").Value = 67.590027700831

$ws.Range("C18").Value = 3
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -21.428571428571
$ws.Range("I18").Value = 181
$ws.Range("J18").Value = 140
$ws.Range("K18").Value = 29.285714285714
$ws.Range("L18").Value = 50.833333333333
$ws.Range("M18").Value = -14.622641509434
$ws.Range("N18").Value = -88.840937114673

$ws.Range("C19").Value = 25
$ws.Range("D19").Value = 16
$ws.Range("E19").Value = 56.25
$ws.Range("F19").Value = 82
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = 6.493506493506
$ws.Range("I19").Value = 807
$ws.Range("J19").Value = 780
$ws.Range("K19").Value = 3.461538461538
$ws.Range("L19").Value = -2.300242130750
$ws.Range("M19").Value = 104.822335025381
$ws.Range("N19").Value = -8.813559322033

$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = -69.230769230769
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 32
$ws.Range("H20").Value = -15.625
$ws.Range("I20").Value = 219
$ws.Range("J20").Value = 246
$ws.Range("K20").Value = -10.975609756097
$ws.Range("L20").Value = 40.384615384615
$ws.Range("M20").Value = 99.090909090909
$ws.Range("N20").Value = -86.506469500924

$ws.Range("C21").Value = 60
$ws.Range("D21").Value = 67
$ws.Range("E21").Value = -10.447761194029
$ws.Range("F21").Value = 221
$ws.Range("G21").Value = 237
$ws.Range("H21").Value = -6.751054852320
$ws.Range("I21").Value = 2227
$ws.Range("J21").Value = 2075
$ws.Range("K21").Value = 7.325301204819
$ws.Range("L21").Value = 22.835079977937
$ws.Range("M21").Value = 88.250211327134
$ws.Range("N21").Value = -60.785349533368

$ws.Range("F22").Value = 3
$ws.Range("H22").Value = -50
$ws.Range("I22").Value = 37
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 37.037037037037
$ws.Range("M22").Value = 54.166666666666

$ws.Range("C24").Value = 49
$ws.Range("D24").Value = 49
$ws.Range("F24").Value = 210
$ws.Range("G24").Value = 229
$ws.Range("H24").Value = -8.296943231441
$ws.Range("I24").Value = 2296
$ws.Range("J24").Value = 2059
$ws.Range("K24").Value = 11.510441962117
$ws.Range("L24").Value = 34.112149532710
$ws.Range("M24").Value = 80.929866036249

$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 36
$ws.Range("E25").Value = -5.555555555555
$ws.Range("F25").Value = 154
$ws.Range("G25").Value = 159
$ws.Range("H25").Value = -3.144654088050
$ws.Range("I25").Value = 1816
$ws.Range("J25").Value = 1502
$ws.Range("K25").Value = 20.905459387483
$ws.Range("L25").Value = 47.642276422764

$ws.Range("C26").Value = 36
$ws.Range("D26").Value = 20
$ws.Range("E26").Value = 80
$ws.Range("F26").Value = 114
$ws.Range("G26").Value = 103
$ws.Range("H26").Value = 10.679611650485
$ws.Range("I26").Value = 1010
$ws.Range("J26").Value = 850
$ws.Range("K26").Value = 18.823529411764
$ws.Range("L26").Value = 54.907975460122
$ws.Range("M26").Value = 116.738197424893

$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -87.5
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = 5.128205128205
$ws.Range("L27").Value = -2.380952380952

$ws.Range("C28").Value = 7
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 600
$ws.Range("F28").Value = 18
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 63.636363636363
$ws.Range("I28").Value = 114
$ws.Range("J28").Value = 102
$ws.Range("K28").Value = 11.764705882352
$ws.Range("L28").Value = 29.545454545454

$ws.Range("L29").Value = -84.615384615384
$ws.Range("N29").Value = -95.833333333333

$ws.Range("L30").Value = -75
$ws.Range("N30").Value = -95.454545454545
